$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PortalTransparencia")

# Update the "Nivel de impacto" (column D) values for rows 31-50 to "Alta"
# and ensure the cells are left-aligned horizontally (matches the style
# used elsewhere in the table for this column).
$range = $ws.Range("D31:D50")
$range.Value = "Alta"
$range.HorizontalAlignment = -4131

# Update the active window view: scroll so row 40 is the top-left visible
# row and select cell D50.
$window = $excel.ActiveWindow
$window.ScrollRow = 40
$ws.Range("D50").Select()
